# TC08_CDS_Filter_PHSAccession-phs002305.xlsx
# commit: "cds SCRIPTS 1-10 obj correction input file correction"
#
# The FilesTab Cypher query (cell B4 on Sheet1) had its
# `experimental_strategies` input filter incorrectly hard-coded to
# ["RNA-Seq"]. Every other query on this sheet (and every other filter
# in this same query) leaves that list empty (= "no restriction"), so
# correct this one back to match: experimental_strategies: [].
#
# Also restore the saved view state: the active/selected cell moves from
# B4 to C4, and the wrapped long-text rows (2-4) are re-pinned to Excel's
# real maximum row height (409.5pt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the experimental_strategies filter in the FilesTab query (B4) ---
$cell = $ws.Range("B4")
$query = $cell.Value2
$badFilter  = 'experimental_strategies: ["RNA-Seq"]'
$fixedFilter = 'experimental_strategies: []'
if ($query.Contains($badFilter)) {
    $cell.Value = $query.Replace($badFilter, $fixedFilter)
}

# --- 2) Rows 2-4 wrap long query text; pin their height at Excel's cap ---
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# --- 3) Move the saved selection from B4 to C4 ---
[void]$ws.Range("C4").Select()
